# Adds a "clarifications needed" column (column N) to the requirements
# table ("Table2") on the Details sheet, with review/clarification notes
# for several of the rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Details")

# The requirements table on this sheet is a structured Table (ListObject)
# spanning A2:M13. Add a 14th column to it -- this naturally extends the
# table ref to A2:N13 and the AutoFilter range along with it.
$lo = $ws.ListObjects.Item(1)
$newCol = $lo.ListColumns.Add()

# Header / title for the new column. Setting the header cell's value is
# what actually renames the ListColumn (it mirrors the table header row).
$ws.Range("N2").Value = "clarifications needed"

# Per-row clarification notes gathered during review.
$ws.Range("N3").Value = "1. what is the use case of locking specific Biometric auths and combinations?`n2. What is the data taken as input from the Individual?`n3. Is there a mechanism to lock OTP Authentication?`n4. need more clarity on a2`n5. Can Lock/Unlock will be perfomed only by OTP authentication of Mobile number or can it be done by email as well?`n6. Need to check with IDA on the process of authentication done, so that the gaps could be covered`n"

$ws.Range("N4").Value = "1. What is the use case of reprinting? Will there be a cost associated to it? If not can it be abused by the individual?will there be a limit on number times an individual can access it? If cost is associated, will there be a check performed for the payment?`n2. Why RID is accepted as an input parameter? What is the use case."

$ws.Range("N5").Value = "1. if demo auth is locked? What happens?`n2. what is use case of providing a RID PDF, why not just a RID number?`n3. why do we have a size check here? Shouldn’t it be stopped at the initial stage?"

$ws.Range("N6").Value = "1. Is this requirement still part of Resident services?`n2. if demo auth is locked? What happens?`n3. why do we have a size check here? Shouldn’t it be stopped at the initial stage?"

$ws.Range("N7").Value = "1. Is this requirement a subset of what registration client does for update? If so can there be an reuse of the feature?`n2. in future if there can be many parameters which can be provisioned for updated? What can be done and is Resident services capable of doing for all? check for scalability?`n"

# Note: N9 is written before N8 so the two short notes land in the shared
# string table in the same relative order as the authored workbook.
$ws.Range("N9").Value = "1.what is security code?"
$ws.Range("N8").Value = "1. why is RID an input here and not UIN?"

# --- Formatting -------------------------------------------------------
# Give the new column a sensible width to hold the multi-line notes.
$ws.Columns.Item(14).ColumnWidth = 47.71

# Match the look of the rest of the table: thin boxed borders, left/top
# aligned text, and wrap the long multi-line notes. Reuse formatting from
# neighbouring cells where possible (copy/paste-special formats) so the
# workbook's style table grows the same way Excel would grow it, instead
# of inventing lots of near-duplicate styles.
$ws.Range("H2").Copy()
$ws.Range("N2").PasteSpecial(-4122)

$ws.Range("N3").WrapText = $true
$ws.Range("N3").HorizontalAlignment = -4131
$ws.Range("N3").VerticalAlignment = -4160
$ws.Range("N3").Borders.LineStyle = 1
$ws.Range("N3").Borders.Weight = 2

$ws.Range("F4").Copy()
$ws.Range("N4:N7").PasteSpecial(-4122)

$ws.Range("H8").Copy()
$ws.Range("N8:N9").PasteSpecial(-4122)

$ws.Range("H10").Copy()
$ws.Range("N10:N12").PasteSpecial(-4122)

$ws.Range("H3").Copy()
$ws.Range("N13").PasteSpecial(-4122)

$excel.CutCopyMode = 0
